$d = $word.ActiveDocument

function Set-ParagraphText($paragraph, $newText) {
    $rng = $paragraph.Range
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

function Replace-ParagraphByOldText($oldText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($ptext -eq $oldText) {
            Set-ParagraphText $p $newText
            return
        }
    }
}

Replace-ParagraphByOldText "LOCADOR:  landlord_name " "LOCADOR:  nome_locador "
Replace-ParagraphByOldText "LOCATÁRIO:  tenant_name " "LOCATÁRIO:  nome_locatario "
Replace-ParagraphByOldText "1. OBJETO: Constitui objeto do presente contrato a locação do imóvel localizado no endereço:  property_address ." "1. OBJETO: Constitui objeto do presente contrato a locação do imóvel localizado no endereço:  endereco_imovel ."
Replace-ParagraphByOldText "2. VALOR DO ALUGUEL: Fica estipulado o valor mensal do aluguel em  rent_amount  ( valor_extenso ) a ser pago pelo LOCATÁRIO ao LOCADOR, mensalmente, até o dia 5 de cada mês." "2. VALOR DO ALUGUEL: Fica estipulado o valor mensal do aluguel em  valor_aluguel  ( valor_aluguel_extenso ) a ser pago pelo LOCATÁRIO ao LOCADOR, mensalmente, até o dia 5 de cada mês."
Replace-ParagraphByOldText "3. PRAZO: O presente contrato terá início em  start_date  e término em  end_date , data em que o LOCATÁRIO deverá desocupar o imóvel e entregá-lo nas mesmas condições em que o recebeu." "3. PRAZO: O presente contrato terá início em  data_inicio  e término em  data_fim , data em que o LOCATÁRIO deverá desocupar o imóvel e entregá-lo nas mesmas condições em que o recebeu."
Replace-ParagraphByOldText "6. FORO: Para dirimir quaisquer controvérsias oriundas do presente contrato, as partes elegem o foro da comarca de  cidade_foro ." "6. FORO: Para dirimir quaisquer controvérsias oriundas do presente contrato, as partes elegem o foro da comarca de  cidade_foro ."
